# Update cryptos list with new price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.818.29"
$ws.Range("E2").Value = "'  +5.34%  "

# Row 3
$ws.Range("D3").Value = "'3.348.54"
$ws.Range("E3").Value = "'  +5.28%  "

# Row 4
$ws.Range("E4").Value = "'  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'570.19"
$ws.Range("E5").Value = "'  +6.95%  "

# Row 6
$ws.Range("D6").Value = "'151.99"
$ws.Range("E6").Value = "'  +5.46%  "

# Row 7
$ws.Range("E7").Value = "'  +0.14%  "

# Row 8
$ws.Range("D8").Value = "'3.351.57"

# Row 9
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "'  -0.53%  "

# Row 10
$ws.Range("D10").Value = "'7.41"
$ws.Range("E10").Value = "'  +1.53%  "

# Row 11
$ws.Range("D11").Value = "'0.117"

# Row 12
$ws.Range("D12").Value = "'0.439"
$ws.Range("E12").Value = "'  +2.42%  "

# Row 13
$ws.Range("D13").Value = "'3.923.57"
$ws.Range("E13").Value = "'  +5.03%  "

# Row 14
$ws.Range("E14").Value = "'  +0.00%  "

# Row 15
$ws.Range("D15").Value = "'26.87"
$ws.Range("E15").Value = "'  +3.94%  "

# Row 16
$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = "'  +4.73%  "

# Row 17
$ws.Range("D17").Value = "'62.839.25"
$ws.Range("E17").Value = "'  +5.25%  "

# Row 18
$ws.Range("D18").Value = "'3.296.18"
$ws.Range("E18").Value = "'  +3.20%  "

# Row 19
$ws.Range("D19").Value = "'6.31"
$ws.Range("E19").Value = "'  +1.77%  "

# Row 20
$ws.Range("D20").Value = "'13.87"
$ws.Range("E20").Value = "'  +5.94%  "

# Row 21
$ws.Range("D21").Value = "'8.38"
$ws.Range("E21").Value = "'  +2.79%  "

# Row 22
$ws.Range("D22").Value = "'384.33"
$ws.Range("E22").Value = "'  +5.11%  "

# Row 23
$ws.Range("D23").Value = "'1.01"
$ws.Range("E23").Value = "'  +0.78%  "

# Row 24
$ws.Range("E24").Value = "'  +2.44%  "

# Row 25
$ws.Range("D25").Value = "'70.39"
$ws.Range("E25").Value = "'  +1.43%  "

# Row 26
$ws.Range("D26").Value = "'9.35"
$ws.Range("E26").Value = "'  +6.85%  "

# Row 27
$ws.Range("E27").Value = "'  +7.00%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0960"
$ws.Range("E28").Value = "'  +9.95%  "

# Row 29
$ws.Range("E29").Value = "'  +0.02%  "

# Row 30
$ws.Range("D30").Value = "'2.00"
$ws.Range("E30").Value = "'  +6.26%  "

# Row 31
$ws.Range("D31").Value = "'22.92"
$ws.Range("E31").Value = "'  +3.26%  "

# Row 32
$ws.Range("D32").Value = "'5.54"
$ws.Range("E32").Value = "'  +5.12%  "

# Row 33
$ws.Range("D33").Value = "'1.30"
$ws.Range("E33").Value = "'  +9.91%  "

# Row 34
$ws.Range("D34").Value = "'6.23"
$ws.Range("E34").Value = "'  +2.79%  "

# Row 35
$ws.Range("D35").Value = "'6.67"
$ws.Range("E35").Value = "'  +2.30%  "

# Row 36
$ws.Range("E36").Value = "'  +9.48%  "

# Row 37
$ws.Range("D37").Value = "'157.43"
$ws.Range("E37").Value = "'  +1.16%  "

# Row 38
$ws.Range("D38").Value = "'1.86"
$ws.Range("E38").Value = "'  +12.22%  "

# Row 39
$ws.Range("D39").Value = "'26.90"

# Row 40
$ws.Range("D40").Value = "'0.0328"
$ws.Range("E40").Value = "'  +12.79%  "

# Row 41
$ws.Range("D41").Value = "'0.0738"
$ws.Range("E41").Value = "'  +6.09%  "

# Row 42
$ws.Range("D42").Value = "'2.816.99"
$ws.Range("E42").Value = "'  +1.26%  "

# Row 43
$ws.Range("D43").Value = "'40.90"
$ws.Range("E43").Value = "'  +3.80%  "

# Row 44
$ws.Range("D44").Value = "'4.25"
$ws.Range("E44").Value = "'  +1.18%  "

# Row 45
$ws.Range("D45").Value = "'0.742"
$ws.Range("E45").Value = "'  +4.61%  "

# Row 46
$ws.Range("B46").Value = "'RenzoRestakedETH"
$ws.Range("C46").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D46").Value = "'3.395.75"
$ws.Range("E46").Value = "'  +5.33%  "

# Row 47
$ws.Range("B47").Value = "'ONDO"
$ws.Range("C47").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "'  +5.37%  "

# Row 48
$ws.Range("D48").Value = "'21.88"
$ws.Range("E48").Value = "'  +7.25%  "

# Row 49
$ws.Range("E49").Value = "'  -2.19%  "

# Row 50
$ws.Range("B50").Value = "'Bittensor"
$ws.Range("C50").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "'293.41"
$ws.Range("E50").Value = "'  +12.12%  "

# Row 51
$ws.Range("B51").Value = "'Cosmos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.29"
$ws.Range("E51").Value = "'  +2.68%  "
